{"js": "// Update the date line and every \"A\u00d7B=C\" multiplication answer in the\n// document to the new values from the commit. Every source string in the\n// document is unique, so an exact, case-sensitive whole-body search for the\n// old value followed by a full replace of that single hit is unambiguous\n// and safe to run against the table cells (and the heading paragraph).\nconst replacements = [\n  ['2025-05-27 Tuesday', '2025-05-28 Wednesday'],\n  ['155\u00d77=1085', '306\u00d72=612'],\n  ['760\u00d75=3800', '285\u00d72=570'],\n  ['563\u00d76=3378', '808\u00d77=5656'],\n  ['315\u00d79=2835', '581\u00d78=4648'],\n  ['437\u00d72=874', '711\u00d77=4977'],\n  ['855\u00d74=3420', '728\u00d75=3640'],\n  ['261\u00d77=1827', '225\u00d79=2025'],\n  ['127\u00d77=889', '743\u00d74=2972'],\n  ['426\u00d76=2556', '952\u00d75=4760'],\n  ['395\u00d79=3555', '196\u00d72=392'],\n  ['954\u00d73=2862', '115\u00d75=575'],\n  ['823\u00d74=3292', '564\u00d77=3948'],\n  ['212\u00d77=1484', '290\u00d76=1740'],\n  ['167\u00d77=1169', '573\u00d78=4584'],\n  ['921\u00d73=2763', '688\u00d73=2064'],\n  ['422\u00d76=2532', '261\u00d79=2349'],\n  ['302\u00d77=2114', '213\u00d77=1491'],\n  ['454\u00d79=4086', '171\u00d77=1197'],\n  ['124\u00d79=1116', '578\u00d78=4624'],\n  ['957\u00d79=8613', '147\u00d77=1029'],\n  ['349\u00d77=2443', '847\u00d78=6776'],\n  ['523\u00d76=3138', '706\u00d74=2824'],\n  ['653\u00d72=1306', '268\u00d72=536'],\n  ['972\u00d73=2916', '521\u00d76=3126'],\n  ['923\u00d72=1846', '394\u00d74=1576'],\n];\n\nconst body = context.document.body;\nconst allResults = [];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  allResults.push({ results, newText });\n}\n\nawait context.sync();\n\nfor (const { results, newText } of allResults) {\n  if (results.items.length === 0) {\n    continue;\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and every \"A\u00d7B=C\" multiplication answer to the new\n# values from the commit. Every old string occurs exactly once in the\n# document, so a plain Find/Replace (Replace:=wdReplaceAll == 2) for each\n# exact, case-sensitive pair touches only the single matching run, whether\n# it's the heading paragraph or a table cell.\n$pairs = @(\n    @(\"2025-05-27 Tuesday\", \"2025-05-28 Wednesday\"),\n    @(\"155\u00d77=1085\", \"306\u00d72=612\"),\n    @(\"760\u00d75=3800\", \"285\u00d72=570\"),\n    @(\"563\u00d76=3378\", \"808\u00d77=5656\"),\n    @(\"315\u00d79=2835\", \"581\u00d78=4648\"),\n    @(\"437\u00d72=874\", \"711\u00d77=4977\"),\n    @(\"855\u00d74=3420\", \"728\u00d75=3640\"),\n    @(\"261\u00d77=1827\", \"225\u00d79=2025\"),\n    @(\"127\u00d77=889\", \"743\u00d74=2972\"),\n    @(\"426\u00d76=2556\", \"952\u00d75=4760\"),\n    @(\"395\u00d79=3555\", \"196\u00d72=392\"),\n    @(\"954\u00d73=2862\", \"115\u00d75=575\"),\n    @(\"823\u00d74=3292\", \"564\u00d77=3948\"),\n    @(\"212\u00d77=1484\", \"290\u00d76=1740\"),\n    @(\"167\u00d77=1169\", \"573\u00d78=4584\"),\n    @(\"921\u00d73=2763\", \"688\u00d73=2064\"),\n    @(\"422\u00d76=2532\", \"261\u00d79=2349\"),\n    @(\"302\u00d77=2114\", \"213\u00d77=1491\"),\n    @(\"454\u00d79=4086\", \"171\u00d77=1197\"),\n    @(\"124\u00d79=1116\", \"578\u00d78=4624\"),\n    @(\"957\u00d79=8613\", \"147\u00d77=1029\"),\n    @(\"349\u00d77=2443\", \"847\u00d78=6776\"),\n    @(\"523\u00d76=3138\", \"706\u00d74=2824\"),\n    @(\"653\u00d72=1306\", \"268\u00d72=536\"),\n    @(\"972\u00d73=2916\", \"521\u00d76=3126\"),\n    @(\"923\u00d72=1846\", \"394\u00d74=1576\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
